$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameter")

# Update the PSO parameter row (Size, Limit, Bloc, Bglob, Brand)
$ws.Range("A2").Value = 20
$ws.Range("B2").Value = 100
$ws.Range("C2").Value = 0.6
$ws.Range("D2").Value = 0.2

# Brand gets one decimal display formatting
$ws.Range("E2").NumberFormat = "0.0"
$ws.Range("E2").Value = 0.2

# Move the selection off the edited row
$ws.Range("E3").Select()
